$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The MidDate recorded for the Essential poll in row 62 (42415, i.e. 24/02/2016)
# was wrong - correct it to 42689 (08/11/2016).
$ws.Range("A62").Value2 = 42689

# With the corrected date this row now belongs further down the (chronologically
# ordered) list. Re-sort the surrounding block of rows by MidDate (column A,
# ascending) so the corrected entry settles into its proper position, shifting
# the intervening rows up by one.
$sortRange = $ws.Range("A61:J70")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A61:A70"), 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

# Update the active selection in the (frozen/split) lower pane to reflect where
# the user was last working.
$ws.Range("F10").Select()
